# Adds bootstrap 95% confidence-interval columns (J:O) to the job-quality
# risk-ratio results sheet: rr_physical_95low, rr_stress_95low,
# rr_poverty_95low, rr_physical_95up, rr_stress_95up, rr_poverty_95up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1), matching the style of existing headers (G1:I1) ---
$headers = @(
    "rr_physical_95low",
    "rr_stress_95low",
    "rr_poverty_95low",
    "rr_physical_95up",
    "rr_stress_95up",
    "rr_poverty_95up"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 10 + $i   # J=10 .. O=15
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter, matches existing header style
}

# --- Data rows 2..19 for columns J..O ---
$data = @{
    2  = @(1.01256515795278,  1.179248466323973, 0.9377009350381036, 1.080764217261617, 1.30371565439118,  1.190528818277899)
    3  = @(1.021397836602722, 1.138869411195226, 0.9954268387171304, 1.094605854342868, 1.253797880451045, 1.217252763449194)
    4  = @(1.016995271843494, 1.160075802067979, 0.9432293953044814, 1.103546737346204, 1.274389571943006, 1.169858813629332)
    5  = @(1.026316720233759, 1.139671311281238, 0.9957410649959413, 1.119083014790386, 1.249695165377421, 1.206099228626335)
    6  = @(1.0173123267935,   1.120546570313376, 0.9303853489482627, 1.110152924532219, 1.211245153329125, 1.217050248481755)
    7  = @(1.027355266201154, 1.098066442082346, 0.9946751112361216, 1.124346328721834, 1.177147881837192, 1.253983410741078)
    8  = @(1.017736679430296, 1.156085726534873, 0.9303805917231224, 1.113031345557702, 1.258940396189938, 1.216384643476736)
    9  = @(1.033218559302173, 1.12352753471734,  0.9948232440375294, 1.14617439113762,  1.216799309851935, 1.24965469312843)
    10 = @(1.020749171424505, 1.136365628382892, 0.9329146872228757, 1.129732230321504, 1.227291662469943, 1.208828122826497)
    11 = @(1.036123910653293, 1.123602155300152, 0.9949415867868214, 1.163601779074461, 1.217905858740534, 1.243724626160237)
    12 = @(1.023403753044698, 1.114250915312967, 0.9270682107036017, 1.145114500369772, 1.188547661986049, 1.228293644953695)
    13 = @(1.040234584110352, 1.096757012981129, 0.9944408472738675, 1.18063482054866,  1.168711053250433, 1.267724618777162)
    14 = @(1.031093895965812, 1.120449468867102, 0.9279574120239251, 1.194587325713757, 1.197876866690254, 1.226603536185932)
    15 = @(1.045447322304079, 1.091864507611691, 0.994423042187409,  1.205024432247845, 1.157843775384524, 1.267317158983811)
    16 = @(1.026774823330058, 1.109460273814597, 0.9312197572538751, 1.174751794305749, 1.194865927990239, 1.214911339040698)
    17 = @(1.040536347952977, 1.098843498333451, 0.9946209566745439, 1.180640130864909, 1.182667975858652, 1.258248855876363)
    18 = @(1.03320098100771,  1.091435160835794, 0.9259802667476098, 1.209345418056715, 1.149304905756631, 1.232243413256689)
    19 = @(1.048790066386728, 1.075379423427375, 0.9943239381972824, 1.223255167240542, 1.12743979681619,  1.273555561072473)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 10 + $i   # J=10 .. O=15
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

Write-Output "Added rr_*_95low / rr_*_95up columns (J:O) for rows 1-19"
